$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-20 Saturday", "2024-04-21 Sunday"),
    @("34×85=", "92×52="),
    @("84×31=", "11×46="),
    @("91×30=", "64×69="),
    @("98×79=", "47×26="),
    @("59×80=", "54×87="),
    @("34×40=", "15×14="),
    @("92×16=", "48×12="),
    @("88×95=", "56×83="),
    @("60×83=", "66×19="),
    @("67×84=", "90×13="),
    @("64×52=", "92×25="),
    @("53×29=", "95×89="),
    @("38×51=", "19×35="),
    @("26×94=", "25×84="),
    @("59×96=", "82×23="),
    @("42×76=", "46×49="),
    @("45×25=", "84×77="),
    @("38×93=", "89×33="),
    @("62×85=", "60×39="),
    @("47×53=", "17×56="),
    @("51×67=", "62×90="),
    @("61×49=", "95×28="),
    @("29×45=", "36×37="),
    @("19×46=", "35×98="),
    @("59×87=", "27×76=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
